$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51 / 52: drop the two stray helper cells P51/P52 ---
$ws.Range("P51").ClearContents()
$ws.Range("P52").ClearContents()

# --- New "Eliminacja <5" label cell (copies the label style used by the
#     neighbouring bucket-name cells, then overwrite with the new text) ---
$ws.Range("C49").Copy($ws.Range("C61"))
$ws.Range("C61").Value = "Eliminacja <5"

# --- Rows 62:66 -- copy down the values of the first five "Oczekiwane"/
#     "Zaobserwowane" buckets (D50:E54) as plain values (not formulas),
#     matching a paste-values style duplication further down the sheet
#     where the low-frequency buckets get consolidated ---
$ws.Range("D62").Value = $ws.Range("D50").Value2
$ws.Range("E62").Value = $ws.Range("E50").Value2
$ws.Range("D63").Value = $ws.Range("D51").Value2
$ws.Range("E63").Value = $ws.Range("E51").Value2
$ws.Range("D64").Value = $ws.Range("D52").Value2
$ws.Range("E64").Value = $ws.Range("E52").Value2
$ws.Range("D65").Value = $ws.Range("D53").Value2
$ws.Range("E65").Value = $ws.Range("E53").Value2
$ws.Range("D66").Value = $ws.Range("D54").Value2
$ws.Range("E66").Value = $ws.Range("E54").Value2

# Row 67 aggregates the remaining low-frequency buckets (rows 55:57)
$ws.Range("D67").Formula = "=SUM(D55:D57)"
$ws.Range("E67").Formula = "=SUM(E55:E57)"

# Chi-square component for each of the new rows
$ws.Range("G62").Formula = "=(POWER(D62-E62,2)/D62)"
$ws.Range("G63:G67").Formula = "=(POWER(D63-E63,2)/D63)"

# --- Row 68: "suma" label + total of the new chi-square components ---
$ws.Range("F58").Copy($ws.Range("F68"))
$ws.Range("G68").Formula = "=SUM(G62:G67)"

# --- View state: move the selection the author ended up on ---
$ws.Range("H73").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
